$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price (D) and Volume(1h) (E) columns to Text format before writing,
# since the source values are plain text (some numeric-looking, e.g. "20.00",
# "0.999") and must not be silently converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "70.112.37"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3
$ws.Range("D3").Value = "3.597.15"
$ws.Range("E3").Value = "  -0.20%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "580.94"
$ws.Range("E5").Value = "  -1.14%  "

# Row 6
$ws.Range("D6").Value = "191.69"
$ws.Range("E6").Value = "  +0.68%  "

# Row 7
$ws.Range("E7").Value = "  -1.41%  "

# Row 8
$ws.Range("D8").Value = "3.594.19"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("E10").Value = "  +3.35%  "

# Row 11
$ws.Range("E11").Value = "  +0.96%  "

# Row 12
$ws.Range("D12").Value = "55.92"
$ws.Range("E12").Value = "  -3.34%  "

# Row 13
$ws.Range("E13").Value = "  +6.63%  "

# Row 14
$ws.Range("D14").Value = "9.70"
$ws.Range("E14").Value = "  -0.74%  "

# Row 15
$ws.Range("D15").Value = "4.180.74"
$ws.Range("E15").Value = "  -0.13%  "

# Row 16
$ws.Range("D16").Value = "20.00"
$ws.Range("E16").Value = "  +3.28%  "

# Row 17
$ws.Range("D17").Value = "3.602.09"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18
$ws.Range("D18").Value = "70.149.64"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19
$ws.Range("E19").Value = "  +2.02%  "

# Row 20
$ws.Range("E20").Value = "  +0.29%  "

# Row 21
$ws.Range("E21").Value = "  -0.11%  "

# Row 22
$ws.Range("D22").Value = "481.51"
$ws.Range("E22").Value = "  -1.87%  "

# Row 23
$ws.Range("D23").Value = "19.23"
$ws.Range("E23").Value = "  +10.91%  "

# Row 24
$ws.Range("D24").Value = "5.04"
$ws.Range("E24").Value = "  -6.43%  "

# Row 25
$ws.Range("E25").Value = "  -0.96%  "

# Row 26
$ws.Range("D26").Value = "95.36"
$ws.Range("E26").Value = "  +5.21%  "

# Row 27
$ws.Range("D27").Value = "3.00"
$ws.Range("E27").Value = "  -3.08%  "

# Row 28
$ws.Range("D28").Value = "11.13"
$ws.Range("E28").Value = "  +0.71%  "

# Row 29
$ws.Range("D29").Value = "9.39"
$ws.Range("E29").Value = "  -0.25%  "

# Row 30
$ws.Range("D30").Value = "32.22"
$ws.Range("E30").Value = "  -0.45%  "

# Row 31
$ws.Range("D31").Value = "7.72"
$ws.Range("E31").Value = "  +3.34%  "

# Row 32
$ws.Range("E32").Value = "  +3.15%  "

# Row 33
$ws.Range("D33").Value = "12.25"

# Row 34
$ws.Range("D34").Value = "66.70"
$ws.Range("E34").Value = "  +2.52%  "

# Row 35
$ws.Range("D35").Value = "587.12"
$ws.Range("E35").Value = "  -6.51%  "

# Row 36
$ws.Range("D36").Value = "39.04"
$ws.Range("E36").Value = "  +3.03%  "

# Row 37
$ws.Range("E37").Value = "  +0.08%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0806"
$ws.Range("E38").Value = "  -1.11%  "

# Row 39
$ws.Range("D39").Value = "0.397"
$ws.Range("E39").Value = "  -1.32%  "

# Row 40
$ws.Range("D40").Value = "3.34"
$ws.Range("E40").Value = "  +23.24%  "

# Row 41
$ws.Range("E41").Value = "  -4.16%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.232.55"
$ws.Range("E42").Value = "  -2.13%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.137"
$ws.Range("E43").Value = "  -6.24%  "

# Row 44
$ws.Range("D44").Value = "2.84"
$ws.Range("E44").Value = "  +6.85%  "

# Row 45
$ws.Range("D45").Value = "3.07"
$ws.Range("E45").Value = "  -0.27%  "

# Row 46
$ws.Range("D46").Value = "0.0449"
$ws.Range("E46").Value = "  +1.39%  "

# Row 47
$ws.Range("D47").Value = "9.49"
$ws.Range("E47").Value = "  +4.22%  "

# Row 48
$ws.Range("E48").Value = "  +0.33%  "

# Row 49
$ws.Range("E49").Value = "  +1.21%  "

# Row 50
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.18%  "

# Row 51
$ws.Range("E51").Value = "  -4.61%  "
